$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 23. This shifts all existing rows
# 23..111 down to 24..112 (and carries the formatting of row 23,
# e.g. the date number format on column D, down with them).
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly record.
# Columns A, B, C, E, F, G, I, R are constant across all rows in this
# sheet (market metadata), so copy them from the (now shifted) row 24.
$ws.Cells.Item(23, 1).Value2 = $ws.Cells.Item(24, 1).Value2   # A: Mercado ID
$ws.Cells.Item(23, 2).Value2 = $ws.Cells.Item(24, 2).Value2   # B: Mercado
$ws.Cells.Item(23, 3).Value2 = $ws.Cells.Item(24, 3).Value2   # C: Región
$ws.Cells.Item(23, 5).Value2 = $ws.Cells.Item(24, 5).Value2   # E: Codreg
$ws.Cells.Item(23, 6).Value2 = $ws.Cells.Item(24, 6).Value2   # F: Categoría ID
$ws.Cells.Item(23, 7).Value2 = $ws.Cells.Item(24, 7).Value2   # G: Categoría
$ws.Cells.Item(23, 9).Value2 = $ws.Cells.Item(24, 9).Value2   # I: Calidad
$ws.Cells.Item(23, 18).Value2 = $ws.Cells.Item(24, 18).Value2 # R: Clasificación

# New/changed values for this record.
$ws.Cells.Item(23, 4).Value2 = 44715                                   # D: Fecha
$ws.Cells.Item(23, 8).Value2 = "Inferno"                               # H: Variedad
$ws.Cells.Item(23, 10).Value2 = 220                                    # J: Volumen
$ws.Cells.Item(23, 11).Value2 = 20000                                  # K: Precio mínimo
$ws.Cells.Item(23, 12).Value2 = 22000                                  # L: Precio máximo
$ws.Cells.Item(23, 13).Value2 = 20909                                  # M: Precio promedio ponderado
$ws.Cells.Item(23, 14).Value2 = "$/caja 12 kilos"                      # N: Unidad de comercialización
$ws.Cells.Item(23, 15).Value2 = "Región de Arica y Parinacota"         # O: Origen
$ws.Cells.Item(23, 16).Value2 = 1742                                   # P: Precio $/Kg
$ws.Cells.Item(23, 17).Value2 = 12                                     # Q: Kg o Unidades
